$d = $word.ActiveDocument
$insertRange = $d.Paragraphs($d.Paragraphs.Count).Range
$insertRange.Collapse(0)

# Pass 1: create paragraphs with plain (non-italic) concatenated text,
# and record absolute start offsets for each paragraph + run boundaries.
$paraInfo = @()

# --- Paragraph 0 ---
$insertRange.InsertParagraphAfter()
$newPara0 = $d.Paragraphs($d.Paragraphs.Count)
$newPara0.Style = "Heading1"
$rng0 = $newPara0.Range
$rng0.Collapse(0)
$rng0.InsertAfter('Knärot – ekologi samt krav på livsmiljön')
$start0 = $newPara0.Range.Start
$insertRange = $newPara0.Range
$insertRange.Collapse(0)

# --- Paragraph 1 ---
$insertRange.InsertParagraphAfter()
$newPara1 = $d.Paragraphs($d.Paragraphs.Count)
$newPara1.Style = "Normal"
$rng1 = $newPara1.Range
$rng1.Collapse(0)
$rng1.InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')
$start1 = $newPara1.Range.Start
$insertRange = $newPara1.Range
$insertRange.Collapse(0)

# --- Paragraph 2 ---
$insertRange.InsertParagraphAfter()
$newPara2 = $d.Paragraphs($d.Paragraphs.Count)
$newPara2.Style = "Normal"
$rng2 = $newPara2.Range
$rng2.Collapse(0)
$rng2.InsertAfter('Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$start2 = $newPara2.Range.Start
$insertRange = $newPara2.Range
$insertRange.Collapse(0)

# --- Paragraph 3 ---
$insertRange.InsertParagraphAfter()
$newPara3 = $d.Paragraphs($d.Paragraphs.Count)
$newPara3.Style = "Normal"
$rng3 = $newPara3.Range
$rng3.Collapse(0)
$rng3.InsertAfter('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$start3 = $newPara3.Range.Start
$insertRange = $newPara3.Range
$insertRange.Collapse(0)

# --- Paragraph 4 ---
$insertRange.InsertParagraphAfter()
$newPara4 = $d.Paragraphs($d.Paragraphs.Count)
$newPara4.Style = "Normal"
$rng4 = $newPara4.Range
$rng4.Collapse(0)
$rng4.InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')
$start4 = $newPara4.Range.Start
$insertRange = $newPara4.Range
$insertRange.Collapse(0)

# --- Paragraph 5 ---
$insertRange.InsertParagraphAfter()
$newPara5 = $d.Paragraphs($d.Paragraphs.Count)
$newPara5.Style = "Normal"
$rng5 = $newPara5.Range
$rng5.Collapse(0)
$rng5.InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')
$start5 = $newPara5.Range.Start
$insertRange = $newPara5.Range
$insertRange.Collapse(0)

# --- Paragraph 6 ---
$insertRange.InsertParagraphAfter()
$newPara6 = $d.Paragraphs($d.Paragraphs.Count)
$newPara6.Style = "Heading2"
$rng6 = $newPara6.Range
$rng6.Collapse(0)
$rng6.InsertAfter('Referenser - knärot')
$start6 = $newPara6.Range.Start
$insertRange = $newPara6.Range
$insertRange.Collapse(0)

# --- Paragraph 7 ---
$insertRange.InsertParagraphAfter()
$newPara7 = $d.Paragraphs($d.Paragraphs.Count)
$newPara7.Style = "Normal"
$rng7 = $newPara7.Range
$rng7.Collapse(0)
$rng7.InsertAfter('de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025')
$start7 = $newPara7.Range.Start
$insertRange = $newPara7.Range
$insertRange.Collapse(0)

# --- Paragraph 8 ---
$insertRange.InsertParagraphAfter()
$newPara8 = $d.Paragraphs($d.Paragraphs.Count)
$newPara8.Style = "Normal"
$rng8 = $newPara8.Range
$rng8.Collapse(0)
$rng8.InsertAfter('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 ')
$start8 = $newPara8.Range.Start
$insertRange = $newPara8.Range
$insertRange.Collapse(0)

# --- Paragraph 9 ---
$insertRange.InsertParagraphAfter()
$newPara9 = $d.Paragraphs($d.Paragraphs.Count)
$newPara9.Style = "Normal"
$rng9 = $newPara9.Range
$rng9.Collapse(0)
$rng9.InsertAfter('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853')
$start9 = $newPara9.Range.Start
$insertRange = $newPara9.Range
$insertRange.Collapse(0)

# --- Paragraph 10 ---
$insertRange.InsertParagraphAfter()
$newPara10 = $d.Paragraphs($d.Paragraphs.Count)
$newPara10.Style = "Normal"
$rng10 = $newPara10.Range
$rng10.Collapse(0)
$rng10.InsertAfter('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.')
$start10 = $newPara10.Range.Start
$insertRange = $newPara10.Range
$insertRange.Collapse(0)

# --- Paragraph 11 ---
$insertRange.InsertParagraphAfter()
$newPara11 = $d.Paragraphs($d.Paragraphs.Count)
$newPara11.Style = "Normal"
$rng11 = $newPara11.Range
$rng11.Collapse(0)
$rng11.InsertAfter('Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$start11 = $newPara11.Range.Start
$insertRange = $newPara11.Range
$insertRange.Collapse(0)

# --- Paragraph 12 ---
$insertRange.InsertParagraphAfter()
$newPara12 = $d.Paragraphs($d.Paragraphs.Count)
$newPara12.Style = "Normal"
$rng12 = $newPara12.Range
$rng12.Collapse(0)
$rng12.InsertAfter('SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala ')
$start12 = $newPara12.Range.Start
$insertRange = $newPara12.Range
$insertRange.Collapse(0)

# Pass 2: apply italic formatting to recorded sub-ranges.

# --- Paragraph 0 italics ---
$pos = $start0
$pos = $pos + 40

# --- Paragraph 1 italics ---
$pos = $start1
$pos = $pos + 684

# --- Paragraph 2 italics ---
$pos = $start2
$pos = $pos + 34
$d.Range($pos, $pos + 82).Font.Italic = $true
$pos = $pos + 82
$pos = $pos + 162
$d.Range($pos, $pos + 205).Font.Italic = $true
$pos = $pos + 205
$pos = $pos + 7
$d.Range($pos, $pos + 118).Font.Italic = $true
$pos = $pos + 118

# --- Paragraph 3 italics ---
$pos = $start3
$pos = $pos + 205
$d.Range($pos, $pos + 865).Font.Italic = $true
$pos = $pos + 865

# --- Paragraph 4 italics ---
$pos = $start4
$pos = $pos + 1337

# --- Paragraph 5 italics ---
$pos = $start5
$pos = $pos + 868

# --- Paragraph 6 italics ---
$pos = $start6
$pos = $pos + 19

# --- Paragraph 7 italics ---
$pos = $start7
$pos = $pos + 33
$d.Range($pos, $pos + 80).Font.Italic = $true
$pos = $pos + 80
$pos = $pos + 44

# --- Paragraph 8 italics ---
$pos = $start8
$pos = $pos + 62
$d.Range($pos, $pos + 114).Font.Italic = $true
$pos = $pos + 114
$pos = $pos + 39

# --- Paragraph 9 italics ---
$pos = $start9
$pos = $pos + 117
$d.Range($pos, $pos + 90).Font.Italic = $true
$pos = $pos + 90
$pos = $pos + 36

# --- Paragraph 10 italics ---
$pos = $start10
$pos = $pos + 54
$d.Range($pos, $pos + 67).Font.Italic = $true
$pos = $pos + 67
$pos = $pos + 38

# --- Paragraph 11 italics ---
$pos = $start11
$pos = $pos + 22
$d.Range($pos, $pos + 35).Font.Italic = $true
$pos = $pos + 35
$pos = $pos + 128

# --- Paragraph 12 italics ---
$pos = $start12
$pos = $pos + 25
$d.Range($pos, $pos + 36).Font.Italic = $true
$pos = $pos + 36
$pos = $pos + 27

# Update date in header (first page header)
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)
    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}
